# Update header label from "Segment Description" to "Segment Name"
# on both the "AppNexus Add_Edit" and "AppNexus Delete" sheets, then
# move the selection/active tab to "AppNexus Delete" (per the diff:
# activeTab moves from sheet index 0 to sheet index 1).

$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AppNexus Add_Edit")
$wsDel = $wb.Worksheets.Item("AppNexus Delete")

$wsAdd.Range("C1").Value = "Segment Name"
$wsDel.Range("B1").Value = "Segment Name"

# Selection on "AppNexus Add_Edit" moves from A3 to C2.
[void]$wsAdd.Range("C2").Select()

# "AppNexus Delete" becomes the active/selected sheet tab.
[void]$wsDel.Activate()
[void]$wsDel.Range("B1").Select()
